$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.593.22'
$ws.Range("E2").Value = '  -0.51%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.510.95'
$ws.Range("E3").Value = '  -2.18%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '619.84'
$ws.Range("E5").Value = '  +2.45%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.24'
$ws.Range("E6").Value = '  -1.18%  '
$ws.Range("E7").Value = '  -1.60%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.506.30'
$ws.Range("E8").Value = '  -2.13%  '
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.197'
$ws.Range("E10").Value = '  -1.74%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.03'
$ws.Range("E11").Value = '  -6.04%  '
$ws.Range("E12").Value = '  -1.86%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '46.29'
$ws.Range("E13").Value = '  -2.27%  '
$ws.Range("E14").Value = '  -1.42%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.081.76'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.38'
$ws.Range("E16").Value = '  -1.16%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '611.08'
$ws.Range("E17").Value = '  -1.71%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.515.81'
$ws.Range("E18").Value = '  -2.07%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '70.639.37'
$ws.Range("E19").Value = '  -0.63%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.121'
$ws.Range("E20").Value = '  +0.96%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.68'
$ws.Range("E21").Value = '  +0.78%  '
$ws.Range("E22").Value = '  -1.26%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.10'
$ws.Range("E23").Value = '  -2.90%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '98.27'
$ws.Range("E24").Value = '  +0.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '15.54'
$ws.Range("E25").Value = '  -4.48%  '
$ws.Range("E26").Value = '  -2.14%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("E28").Value = '  -3.88%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.65'
$ws.Range("E29").Value = '  -1.96%  '
$ws.Range("E30").Value = '  -3.56%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.02'
$ws.Range("E31").Value = '  -2.36%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.07'
$ws.Range("E32").Value = '  -5.88%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.30'
$ws.Range("E33").Value = '  -0.86%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '638.35'
$ws.Range("E34").Value = '  +1.43%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.77'
$ws.Range("E35").Value = '  -6.30%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0996'
$ws.Range("E36").Value = '  -2.97%  '
$ws.Range("E37").Value = '  -1.28%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0474'
$ws.Range("E38").Value = '  -3.05%  '
$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '56.77'
$ws.Range("E39").Value = '  -1.28%  '
$ws.Range("B40").Value = 'dogwifhat'
$ws.Range("C40").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.41'
$ws.Range("E40").Value = '  -11.22%  '
$ws.Range("E41").Value = '  +0.22%  '
$ws.Range("E42").Value = '  +0.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.339.88'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0₃0716'
$ws.Range("E44").Value = '  -0.62%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.95'
$ws.Range("E45").Value = '  -2.82%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.312'
$ws.Range("E46").Value = '  -4.65%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '31.76'
$ws.Range("E47").Value = '  -4.28%  '
$ws.Range("E48").Value = '  -7.16%  '
$ws.Range("E49").Value = '  -1.09%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '133.70'
$ws.Range("E50").Value = '  +0.66%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.155'
$ws.Range("E51").Value = '  +4.81%  '
